$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header text updates (row 1)
$ws.Range("B1").Value = 'Government-Cadw'
$ws.Range("C1").Value = 'Government-Local_Authority'
$ws.Range("D1").Value = 'Government-National'
$ws.Range("E1").Value = 'Government-Other'
$ws.Range("F1").Value = 'Independent-English_Heritage'
$ws.Range("G1").Value = 'Independent-Historic_Environment_Scotland'
$ws.Range("H1").Value = 'Independent-National_Trust'
$ws.Range("I1").Value = 'Independent-National_Trust_for_Scotland'
$ws.Range("J1").Value = 'Independent-Not_for_profit'
$ws.Range("K1").Value = 'Independent-Private'
$ws.Range("L1").Value = 'Independent-Unknown'
$ws.Range("M1").Value = 'University'
$ws.Range("N1").Value = 'Unknown'

# Numeric value updates (rows 2-6)
$ws.Range("C2").Value = 15.718
$ws.Range("D2").Value = 1.34
$ws.Range("E2").Value = 0.167
$ws.Range("F2").Value = 1.268
$ws.Range("H2").Value = 3.923
$ws.Range("J2").Value = 32.081
$ws.Range("K2").Value = 14.067
$ws.Range("L2").Value = 3.947
$ws.Range("M2").Value = 1.818
$ws.Range("N2").Value = 1.986
$ws.Range("O2").Value = 76.315
$ws.Range("C3").Value = 0.8129999999999999
$ws.Range("D3").Value = 0.096
$ws.Range("H3").Value = 0.167
$ws.Range("J3").Value = 0.766
$ws.Range("K3").Value = 0.431
$ws.Range("L3").Value = 0.07199999999999999
$ws.Range("M3").Value = 0.024
$ws.Range("N3").Value = 0.191
$ws.Range("O3").Value = 2.56
$ws.Range("C4").Value = 3.876
$ws.Range("D4").Value = 0.263
$ws.Range("G4").Value = 0.502
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.646
$ws.Range("J4").Value = 6.005
$ws.Range("K4").Value = 1.842
$ws.Range("L4").Value = 0.9330000000000001
$ws.Range("M4").Value = 0.67
$ws.Range("N4").Value = 0.167
$ws.Range("O4").Value = 14.904
$ws.Range("B5").Value = 0.07199999999999999
$ws.Range("C5").Value = 1.459
$ws.Range("D5").Value = 0.263
$ws.Range("E5").Value = 0.024
$ws.Range("H5").Value = 0.335
$ws.Range("J5").Value = 2.153
$ws.Range("K5").Value = 1.34
$ws.Range("L5").Value = 0.191
$ws.Range("M5").Value = 0.12
$ws.Range("N5").Value = 0.263
$ws.Range("O5").Value = 6.220000000000001
$ws.Range("B6").Value = 0.07199999999999999
$ws.Range("C6").Value = 21.866
$ws.Range("D6").Value = 1.962
$ws.Range("E6").Value = 0.191
$ws.Range("F6").Value = 1.268
$ws.Range("G6").Value = 0.502
$ws.Range("H6").Value = 4.425
$ws.Range("I6").Value = 0.646
$ws.Range("J6").Value = 41.005
$ws.Range("K6").Value = 17.68
$ws.Range("L6").Value = 5.143
$ws.Range("M6").Value = 2.632
$ws.Range("N6").Value = 2.607
$ws.Range("O6").Value = 99.999
